$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3632.353
$ws.Range("J17").Value = 3796.875
$ws.Range("L17").Value = 11390.625
$ws.Range("N17").Value = -11726.625

$ws.Range("H62").Value = 38101452
$ws.Range("I62").Value = 38101452
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 38101452
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = $null
$ws.Range("N62").Value = -38100828

$ws.Range("H65").Value = 38101452
$ws.Range("I65").Value = 38101452
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 190507260
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -190504140
$ws.Range("N65").Value = 0

$ws.Range("H132").Value = 1662
$ws.Range("J132").Value = 1005
$ws.Range("L132").Value = 3015
$ws.Range("N132").Value = -8075

$ws.Range("H138").Value = 1920.9025
$ws.Range("J138").Value = 3328.4546
$ws.Range("L138").Value = 9985.363799999999
$ws.Range("N138").Value = -20265.3638

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H12").Value = 3312.4
$ws.Range("I12").Value = 1779
$ws.Range("J12").Value = 4334.6665
$ws.Range("K12").Value = 1779
$ws.Range("L12").Value = 4334.6665
$ws.Range("M12").Value = -1606
$ws.Range("N12").Value = -4680.6665

$ws.Range("H32").Value = 41563.41
$ws.Range("J32").Value = 17999.8
$ws.Range("L32").Value = 17999.8
$ws.Range("N32").Value = -18573.8

$ws.Range("H35").Value = 1679
$ws.Range("I35").Value = 1679
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 1679
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = $null
$ws.Range("N35").Value = -1273

$ws.Range("H61").Value = 12828022
$ws.Range("I61").Value = 20838462
$ws.Range("J61").Value = 11317
$ws.Range("K61").Value = 20838462
$ws.Range("L61").Value = 11317
$ws.Range("M61").Value = -20838250
$ws.Range("N61").Value = -11741

$ws.Range("H74").Value = 4945.4116
$ws.Range("I74").Value = 1138.2
$ws.Range("K74").Value = 1138.2
$ws.Range("M74").Value = -264.2

$ws.Range("H77").Value = 4945.4116
$ws.Range("I77").Value = 1138.2
$ws.Range("K77").Value = 5691
$ws.Range("M77").Value = -1323

$ws.Range("H97").Value = 2850016.2
$ws.Range("J97").Value = 1350
$ws.Range("L97").Value = 1350
$ws.Range("N97").Value = -2342

$ws.Range("H112").Value = 39500
$ws.Range("J112").Value = 39500
$ws.Range("L112").Value = 39500
$ws.Range("N112").Value = -42454

$ws.Range("H132").Value = 3778011.8
$ws.Range("I132").Value = 5131423.5
$ws.Range("J132").Value = 7792.357
$ws.Range("K132").Value = 15394270.5
$ws.Range("L132").Value = 23377.071
$ws.Range("M132").Value = -15391740.5
$ws.Range("N132").Value = -28437.071

$ws.Range("H136").Value = 12828022
$ws.Range("I136").Value = 20838462
$ws.Range("J136").Value = 11317
$ws.Range("K136").Value = 62515386
$ws.Range("L136").Value = 33951
$ws.Range("M136").Value = -62512836
$ws.Range("N136").Value = -39051

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1569.4828
$ws.Range("I107").Value = 1248.5714
$ws.Range("K107").Value = 1248.5714
$ws.Range("M107").Value = 671.4286

$ws.Range("H123").Value = 85916.664
$ws.Range("J123").Value = 85916.664
$ws.Range("L123").Value = 85916.664
$ws.Range("N123").Value = -95716.664

$ws.Range("H124").Value = 47780
$ws.Range("J124").Value = 47780
$ws.Range("L124").Value = 47780
$ws.Range("N124").Value = -57600

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = $null
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = 0

$ws.Range("H26").Value = 0
$ws.Range("J26").Value = 0
$ws.Range("L26").Value = $null
$ws.Range("N26").Value = 0

$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = $null
$ws.Range("M27").Value = $null
$ws.Range("N27").Value = 0

$ws.Range("H29").Value = 13500
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 13500
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = $null
$ws.Range("M29").Value = 13500
$ws.Range("N29").Value = -14086

$ws.Range("H33").Value = 4000
$ws.Range("I33").Value = 4000
$ws.Range("K33").Value = 4000
$ws.Range("M33").Value = -3621

$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("L36").Value = $null
$ws.Range("M36").Value = $null
$ws.Range("N36").Value = 0

$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = $null
$ws.Range("M40").Value = $null
$ws.Range("N40").Value = 0

$ws.Range("H51").Value = 78346.664
$ws.Range("J51").Value = 94997.5
$ws.Range("L51").Value = 94997.5
$ws.Range("N51").Value = -96469.5

$ws.Range("H61").Value = 78346.664
$ws.Range("J61").Value = 94997.5
$ws.Range("L61").Value = 94997.5
$ws.Range("N61").Value = -95693.5

$ws.Range("H86").Value = 6287.8335
$ws.Range("I86").Value = 3966.6667
$ws.Range("J86").Value = 7061.5557
$ws.Range("K86").Value = 3966.6667
$ws.Range("L86").Value = 7061.5557
$ws.Range("M86").Value = -2843.6667
$ws.Range("N86").Value = -9307.555700000001

$ws.Range("H89").Value = 6287.8335
$ws.Range("I89").Value = 3966.6667
$ws.Range("J89").Value = 7061.5557
$ws.Range("K89").Value = 19833.3335
$ws.Range("L89").Value = 35307.7785
$ws.Range("M89").Value = -14217.3335
$ws.Range("N89").Value = -46539.7785

$ws.Range("H100").Value = 55500
$ws.Range("J100").Value = 55500
$ws.Range("L100").Value = 55500
$ws.Range("N100").Value = -57664

$ws.Range("H130").Value = 60142.855
$ws.Range("J130").Value = 60142.855
$ws.Range("L130").Value = 60142.855
$ws.Range("N130").Value = -70182.85500000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 1702
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").Value = $null

$ws.Range("H83").Value = 1702
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").Value = $null

$ws.Range("H98").Value = 517.1429000000001
$ws.Range("J98").Value = 517.1429000000001
$ws.Range("L98").Value = 1551.4287
$ws.Range("N98").Value = -4547.4287

$ws.Range("H115").Value = 948.8
$ws.Range("I115").Value = 948.8
$ws.Range("K115").Value = 2846.4
$ws.Range("M115").Value = -1671.4

$ws.Range("H129").Value = 29412954
$ws.Range("J129").Value = 50001684
$ws.Range("L129").Value = 150005052
$ws.Range("N129").Value = -150015052

$ws.Range("H132").Value = 56126.684
$ws.Range("I132").Value = 86871.5
$ws.Range("K132").Value = 781843.5
$ws.Range("M132").Value = -779313.5

$ws.Range("H140").Value = 1857.125
$ws.Range("J140").Value = 1963.25
$ws.Range("L140").Value = 5889.75
$ws.Range("N140").Value = -16249.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2429.9707
$ws.Range("I122").Value = 1943.52
$ws.Range("J122").Value = 3781.2222
$ws.Range("K122").Value = 5830.559999999999
$ws.Range("L122").Value = 11343.6666
$ws.Range("M122").Value = -3380.559999999999
$ws.Range("N122").Value = -16243.6666

$ws.Range("H134").Value = 69299.14
$ws.Range("J134").Value = 69299.14
$ws.Range("L134").Value = 207897.42
$ws.Range("N134").Value = -212967.42

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H26").Value = 25000
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 25000
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = $null
$ws.Range("M26").Value = 25000
$ws.Range("N26").Value = -25590

$ws.Range("H56").Value = 24490.5
$ws.Range("I56").Value = 8981
$ws.Range("J56").Value = 40000
$ws.Range("K56").Value = 8981
$ws.Range("L56").Value = 40000
$ws.Range("M56").Value = -8290
$ws.Range("N56").Value = -41382

$ws.Range("H93").Value = 2042.1428
$ws.Range("J93").Value = 1749
$ws.Range("L93").Value = 1749
$ws.Range("N93").Value = -4245

$ws.Range("H124").Value = 76292.664
$ws.Range("J124").Value = 76292.664
$ws.Range("L124").Value = 76292.664
$ws.Range("N124").Value = -86112.664

$ws.Range("H132").Value = 4699.2446
$ws.Range("I132").Value = 4142.5713
$ws.Range("K132").Value = 12427.7139
$ws.Range("M132").Value = -9897.713899999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 4261.5
$ws.Range("J41").Value = 4261.5
$ws.Range("L41").Value = 4261.5
$ws.Range("N41").Value = -5041.5

$ws.Range("H47").Value = 18666.666
$ws.Range("J47").Value = 18666.666
$ws.Range("L47").Value = 18666.666
$ws.Range("N47").Value = -19810.666

$ws.Range("H120").Value = 50960
$ws.Range("J120").Value = 50960
$ws.Range("L120").Value = 50960
$ws.Range("N120").Value = -60636

$ws.Range("H122").Value = 2686.7742
$ws.Range("I122").Value = 2734.476
$ws.Range("J122").Value = 2586.6
$ws.Range("K122").Value = 8203.428
$ws.Range("L122").Value = 7759.799999999999
$ws.Range("M122").Value = -5753.428
$ws.Range("N122").Value = -12659.8

$ws.Range("H132").Value = 6607.439
$ws.Range("I132").Value = 6025.1763
$ws.Range("J132").Value = 9435.571
$ws.Range("K132").Value = 18075.5289
$ws.Range("L132").Value = 28306.713
$ws.Range("M132").Value = -15545.5289
$ws.Range("N132").Value = -33366.713

$ws.Range("H136").Value = 3486550.8
$ws.Range("I136").Value = 4927235
$ws.Range("K136").Value = 14781705
$ws.Range("M136").Value = -14779155
